$d = $word.ActiveDocument

# Locate the "Screen Shots" paragraph and collapse the range to its end
# (i.e. right after that paragraph, before whatever follows it).
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("Screen Shots", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$findRange.Collapse(0)

# Insert a new 5-row x 2-column table right there, matching the table
# already used elsewhere in this document (TableGrid style / 04A0 look).
$tbl = $d.Tables.Add($findRange, 5, 2)
$tbl.Style = "TableGrid"

# Reproduce Word's default "04A0" table look (header row + first column
# emphasis, horizontal banding only) on the new table.
$tbl.ApplyStyleHeadingRows = $true
$tbl.ApplyStyleFirstColumn = $true
$tbl.ApplyStyleLastRow = $false
$tbl.ApplyStyleLastColumn = $false
$tbl.ApplyStyleRowBands = $true
$tbl.ApplyStyleColumnBands = $false

# Match the column widths used by the existing table in the document
# (6 columns x 1596 twips = 9576 twips total, split across 2 columns here).
$tbl.Columns(1).Width = 239.4
$tbl.Columns(2).Width = 239.4

# Fill in the header row and the three field rows; the final row is left blank.
$tbl.Cell(1, 1).Range.Text = "Sno"
$tbl.Cell(1, 2).Range.Text = "Field Name"
$tbl.Cell(2, 1).Range.Text = "1"
$tbl.Cell(2, 2).Range.Text = "Name"
$tbl.Cell(3, 1).Range.Text = "2"
$tbl.Cell(3, 2).Range.Text = "Password"
$tbl.Cell(4, 1).Range.Text = "3"
$tbl.Cell(4, 2).Range.Text = "Email"

# Bold every cell (and the trailing paragraph marks) to match the rest of
# the document's formatting.
$tbl.Range.Font.Bold = 1
